$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.386567831039429
$ws.Range("B1").Value = 2.731710910797119
$ws.Range("C1").Value = 5.641032695770264
$ws.Range("D1").Value = 2.176080226898193
$ws.Range("E1").Value = 1.204565286636353
